$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("C1").Value = "rules"
$ws.Range("D1").Value = "adaptive_filter"

# Column D (rows 2-8): numeric 2 -> text "wRLS"
$ws.Range("D2").Value = "wRLS"
$ws.Range("D3").Value = "wRLS"
$ws.Range("D4").Value = "wRLS"
$ws.Range("D5").Value = "wRLS"
$ws.Range("D6").Value = "wRLS"
$ws.Range("D7").Value = "wRLS"
$ws.Range("D8").Value = "wRLS"

# Updated RMSE / NDEI / MAE values
$ws.Range("E2").Value = 0.3389426223923964
$ws.Range("F2").Value = 0.9462728344299011
$ws.Range("G2").Value = 0.2662458878999209

$ws.Range("E3").Value = 0.3450909056852412
$ws.Range("F3").Value = 0.9634378442989248
$ws.Range("G3").Value = 0.2507222847087046

$ws.Range("E4").Value = 0.3435806533029622
$ws.Range("F4").Value = 0.9592214645695286
$ws.Range("G4").Value = 0.2494468821090003

$ws.Range("E5").Value = 0.341235608020874
$ws.Range("F5").Value = 0.9526744784446057
$ws.Range("G5").Value = 0.2435756649949328

$ws.Range("E6").Value = 0.3543855117192541
$ws.Range("F6").Value = 0.9893868770131768
$ws.Range("G6").Value = 0.2535839278984292

$ws.Range("E7").Value = 0.3397578123769143
$ws.Range("F7").Value = 0.9485487126649954
$ws.Range("G7").Value = 0.2431092099085204

$ws.Range("E8").Value = 0.3405361950855589
$ws.Range("F8").Value = 0.950721830954992
$ws.Range("G8").Value = 0.2658466593202523
